$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content is safe to assign directly (Excel will not
# mis-parse them as a number/date), covering Coin/Link swaps and all
# Volume(1h) percentage cells plus the few Price cells that contain more
# than one '.' (so Excel keeps them as text automatically).
$directValues = [ordered]@{
    'D2' = '28.185.87'
    'E2' = '  +0.19%  '
    'D3' = '1.869.51'
    'E3' = '  +3.04%  '
    'E4' = '  +0.08%  '
    'E6' = '  +0.25%  '
    'E7' = '  -0.96%  '
    'E8' = '  +0.43%  '
    'E9' = '  +0.30%  '
    'E10' = '  +3.10%  '
    'E11' = '  -0.26%  '
    'E12' = '  +0.89%  '
    'E13' = '  +2.50%  '
    'D14' = '1.873.44'
    'E14' = '  +3.15%  '
    'E15' = '  +0.06%  '
    'E16' = '  +1.26%  '
    'E17' = '  -0.21%  '
    'E18' = '  +0.60%  '
    'E19' = '  +0.46%  '
    'E20' = '  +1.62%  '
    'E21' = '  +0.06%  '
    'D23' = '28.244.34'
    'E23' = '  +0.35%  '
    'E24' = '  +1.83%  '
    'E25' = '  +2.75%  '
    'E26' = '  +3.92%  '
    'D27' = '2.084.65'
    'E27' = '  +3.39%  '
    'E28' = '  +3.19%  '
    'E29' = '  +0.36%  '
    'E31' = '  +1.65%  '
    'E32' = '  -3.24%  '
    'E33' = '  -0.12%  '
    'E34' = '  +0.09%  '
    'E35' = '  +5.75%  '
    'E36' = '  -1.92%  '
    'E37' = '  +2.48%  '
    'E38' = '  +0.60%  '
    'B39' = 'TheSandbox'
    'C39' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'E39' = '  +3.01%  '
    'B40' = 'Aptos'
    'C40' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'E40' = '  -0.49%  '
    'E41' = '  -1.03%  '
    'E42' = '  +2.75%  '
    'E43' = '  +0.21%  '
    'E44' = '  +2.35%  '
    'E45' = '  +1.35%  '
    'E46' = '  -0.83%  '
    'E47' = '  -1.61%  '
    'E48' = '  -0.72%  '
    'E49' = '  +1.92%  '
    'E50' = '  +1.16%  '
    'E51' = '  +1.03%  '
}

foreach ($cellRef in $directValues.Keys) {
    $ws.Range($cellRef).Value = $directValues[$cellRef]
}

# Price cells whose new text looks like a plain number (e.g. '311.90').
# The source workbook stores Price as literal text, so force the Text
# number format before assigning, otherwise Excel auto-converts the
# string into a numeric value and the digits after assignment would
# drift (e.g. floating point noise) from the literal text we want.
$textPriceValues = [ordered]@{
    'D4' = '1.001'
    'D5' = '311.90'
    'D8' = '0.3918'
    'D9' = '0.09665'
    'D10' = '1.137'
    'D11' = '40.82'
    'D12' = '6.501'
    'D13' = '20.98'
    'D15' = '1.001'
    'D16' = '7.421'
    'D18' = '92.98'
    'D19' = '0.06628'
    'D20' = '17.54'
    'D22' = '6.154'
    'D24' = '11.33'
    'D25' = '2.281'
    'D26' = '2.528'
    'D28' = '21.20'
    'D29' = '157.60'
    'D30' = '127.44'
    'D31' = '1.071'
    'D32' = '0.1055'
    'D33' = '5.627'
    'D34' = '3.628'
    'D35' = '9.593'
    'D36' = '0.06750'
    'D37' = '0.02388'
    'D38' = '0.2182'
    'D39' = '0.6361'
    'D40' = '11.46'
    'D41' = '4.977'
    'D42' = '1.182'
    'D43' = '1.002'
    'D44' = '13.62'
    'D45' = '0.6021'
    'D47' = '1.266'
    'D48' = '124.12'
    'D49' = '1.994'
    'D51' = '0.06837'
}

foreach ($cellRef in $textPriceValues.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $textPriceValues[$cellRef]
}
